# Daily attendance processing - reorder "Recorded By" entries in column G.
#
# For every populated cell in column G ("Recorded By") whose value is a
# comma-separated list that includes a "System"/"system" entry alongside at
# least one other entry, the list order is reversed (so the "System" marker
# that used to lead now trails the human/automation accounts).
# Single-value cells (e.g. just "System") and multi-value cells that do not
# mention "System" at all are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $firstRow + $usedRange.Rows.Count - 1

# Column G is "Recorded By" (1=A ... 7=G); skip the header row.
$col = 7

for ($r = [Math]::Max($firstRow, 2); $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $text = $cell.Text

    if ($text -eq $null -or $text -eq "") {
        continue
    }

    $rawParts = $text -split ","
    if ($rawParts.Count -le 1) {
        continue
    }

    $parts = @()
    foreach ($p in $rawParts) {
        $parts += $p.Trim()
    }

    $hasSystem = $false
    foreach ($p in $parts) {
        if ($p.ToLower() -eq "system") {
            $hasSystem = $true
        }
    }

    if (-not $hasSystem) {
        continue
    }

    $reversed = @()
    for ($i = $parts.Count - 1; $i -ge 0; $i--) {
        $reversed += $parts[$i]
    }

    $newValue = [string]::Join(", ", $reversed)
    if ($newValue -ne $text) {
        $cell.Value = $newValue
    }
}
